$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 28 (2025Q2) metrics to reflect refreshed Bibi/Add data
$ws.Range("C28").Value = 355
$ws.Range("D28").Value = 38
$ws.Range("E28").Value = 317
$ws.Range("F28").Value = 5.919003115264798
